$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = '[''album'', ''song'', ''madonna'', ''chart'', ''video'', ''music'', ''band'', ''songs'', ''harrison'', ''carey'', ''track'', ''pop'', ''recording'', ''vocals'', ''lyrics'']'
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1

$ws.Cells.Item(3, 2).Value = '[''ship'', ''guns'', ''ships'', ''tons'', ''torpedo'', ''knots'', ''inch'', ''cruiser'', ''fleet'', ''gun'', ''deck'', ''admiral'', ''cruisers'', ''turrets'', ''german'']'
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1

$ws.Cells.Item(4, 2).Value = '[''episode'', ''mulder'', ''scully'', ''doctor'', ''episodes'', ''trek'', ''series'', ''enterprise'', ''character'', ''viewers'', ''files'', ''amy'', ''television'', ''fringe'', ''scene'']'
$ws.Cells.Item(4, 3).Value = 0.4953602701949368
$ws.Cells.Item(4, 4).Value = 0.4929879911887854
$ws.Cells.Item(4, 5).Value = 0.492386805563901
$ws.Cells.Item(4, 6).Value = 0.4871538194895848
$ws.Cells.Item(4, 7).Value = 0.4871538194895848

$ws.Cells.Item(5, 2).Value = '[''game'', ''player'', ''gameplay'', ''games'', ''players'', ''playstation'', ''nintendo'', ''released'', ''graphics'', ''characters'', ''soundtrack'', ''xbox'', ''mode'', ''version'', ''manga'']'
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0.5043347996538532
$ws.Cells.Item(5, 6).Value = 0.296402122081634
$ws.Cells.Item(5, 7).Value = 0.2331156308035778

$ws.Cells.Item(6, 2).Value = '[''innings'', ''team'', ''runs'', ''match'', ''league'', ''nba'', ''baseball'', ''batting'', ''season'', ''career'', ''scored'', ''wickets'', ''championship'', ''cricket'', ''basketball'']'
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 0.8868632889016513
$ws.Cells.Item(6, 6).Value = 0.530041659289464
$ws.Cells.Item(6, 7).Value = 0.5190010126251009

$ws.Cells.Item(7, 2).Value = '[''highway'', ''route'', ''road'', ''freeway'', ''interchange'', ''intersection'', ''terminus'', ''north'', ''east'', ''continues'', ''lane'', ''state'', ''avenue'', ''passes'', ''crosses'']'
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1

$ws.Cells.Item(8, 2).Value = '[''tropical'', ''storm'', ''hurricane'', ''winds'', ''depression'', ''cyclone'', ''mph'', ''rainfall'', ''damage'', ''landfall'', ''utc'', ''wind'', ''flooding'', ''weakened'', ''intensity'']'
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1

$ws.Cells.Item(9, 2).Value = '[''film'', ''films'', ''disney'', ''movie'', ''animated'', ''role'', ''comedy'', ''simpsons'', ''starred'', ''cast'', ''character'', ''production'', ''script'', ''actor'', ''director'']'
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0

$ws.Cells.Item(10, 2).Value = '[''species'', ''shark'', ''genus'', ''sharks'', ''cap'', ''fruit'', ''stem'', ''fungus'', ''spores'', ''brown'', ''nest'', ''females'', ''habitat'', ''eggs'', ''fin'']'
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1

$ws.Cells.Item(11, 2).Value = '[''election'', ''campaign'', ''political'', ''party'', ''bush'', ''government'', ''republican'', ''vote'', ''presidential'', ''president'', ''labour'', ''democratic'', ''senate'', ''candidate'', ''women'']'
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0

$ws.Cells.Item(12, 2).Value = '[''episode'', ''dwight'', ''glee'', ''michael'', ''jim'', ''liz'', ''pam'', ''nbc'', ''viewers'', ''andy'', ''jack'', ''office'', ''watched'', ''jenna'', ''fey'']'
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1

$ws.Cells.Item(13, 2).Value = '[''king'', ''polish'', ''arab'', ''son'', ''emperor'', ''reign'', ''pope'', ''poland'', ''henry'', ''hungary'', ''royal'', ''constantine'', ''byzantine'', ''bishop'', ''died'']'
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0

$ws.Cells.Item(14, 2).Value = '[''french'', ''persian'', ''army'', ''british'', ''battle'', ''troops'', ''fleet'', ''siege'', ''militia'', ''force'', ''men'', ''cavalry'', ''forces'', ''expedition'', ''ships'']'
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1

$ws.Cells.Item(15, 2).Value = '[''hitler'', ''commander'', ''officer'', ''war'', ''physics'', ''command'', ''holocaust'', ''german'', ''promoted'', ''nuclear'', ''lieutenant'', ''jews'', ''awarded'', ''nazi'', ''germany'']'
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0

$ws.Cells.Item(16, 2).Value = '[''building'', ''mosque'', ''museum'', ''temple'', ''pier'', ''square'', ''chicago'', ''park'', ''hall'', ''library'', ''mall'', ''memorial'', ''buildings'', ''memorials'', ''floor'']'
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = 0.1848809736710457
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 0

$ws.Cells.Item(17, 2).Value = '[''castle'', ''century'', ''church'', ''castles'', ''scotland'', ''tower'', ''bailey'', ''stone'', ''walls'', ''scottish'', ''built'', ''painting'', ''chancel'', ''wall'', ''medieval'']'
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 1

$ws.Cells.Item(18, 2).Value = '[''creek'', ''watershed'', ''river'', ''dam'', ''bridge'', ''flows'', ''volcano'', ''area'', ''lake'', ''park'', ''feet'', ''water'', ''canal'', ''city'', ''lava'']'
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.583906551538591
$ws.Cells.Item(18, 7).Value = 0.4981189741592873

$ws.Cells.Item(19, 2).Value = '[''battalion'', ''brigade'', ''division'', ''regiment'', ''infantry'', ''battalions'', ''units'', ''unit'', ''training'', ''squadron'', ''1st'', ''2nd'', ''artillery'', ''casualties'', ''forces'']'
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 1

$ws.Cells.Item(20, 2).Value = '[''film'', ''tamil'', ''telugu'', ''chinese'', ''india'', ''han'', ''hindu'', ''indian'', ''films'', ''mumbai'', ''hindi'', ''china'', ''cinema'', ''sanskrit'', ''bollywood'']'
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 0

$ws.Cells.Item(21, 2).Value = '[''station'', ''trains'', ''railway'', ''line'', ''locomotives'', ''services'', ''train'', ''locomotive'', ''oslo'', ''railways'', ''passenger'', ''platforms'', ''tunnel'', ''nok'', ''platform'']'
$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 1

$ws.Cells.Item(22, 2).Value = '[''oxford'', ''cambridge'', ''race'', ''boat'', ''blues'', ''rowed'', ''rowing'', ''lengths'', ''crews'', ''rower'', ''thames'', ''races'', ''crew'', ''universities'', ''umpired'']'
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(22, 5).Value = 1
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 1

$ws.Cells.Item(23, 2).Value = '[''aircraft'', ''engine'', ''flight'', ''spacecraft'', ''wing'', ''fuselage'', ''mission'', ''fuel'', ''engines'', ''prototype'', ''air'', ''nasa'', ''radar'', ''apollo'', ''landing'']'
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 1

$ws.Cells.Item(24, 2).Value = '[''breed'', ''horses'', ''horse'', ''breeds'', ''dog'', ''dogs'', ''stud'', ''breeding'', ''bred'', ''arabian'', ''stallion'', ''breeders'', ''riding'', ''stakes'', ''pony'']'
$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 1

$ws.Cells.Item(25, 2).Value = '[''governor'', ''kentucky'', ''massachusetts'', ''colony'', ''boston'', ''virginia'', ''elected'', ''fraternity'', ''plymouth'', ''davis'', ''colonial'', ''served'', ''legislature'', ''confederate'', ''state'']'
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 5).Value = 1
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 1

$ws.Cells.Item(26, 2).Value = '[''lap'', ''race'', ''drivers'', ''laps'', ''pit'', ''driver'', ''car'', ''ferrari'', ''qualifying'', ''prix'', ''session'', ''fastest'', ''ahead'', ''hamilton'', ''caution'']'
$ws.Cells.Item(26, 3).Value = 1
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 5).Value = 1
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 1

$ws.Cells.Item(27, 2).Value = '[''poem'', ''poems'', ''poetry'', ''poet'', ''shakespeare'', ''ode'', ''riley'', ''smart'', ''narrator'', ''published'', ''poetic'', ''stanza'', ''literary'', ''thomas'', ''works'']'
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = 1
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 7).Value = 1

$ws.Cells.Item(28, 2).Value = '[''match'', ''event'', ''ring'', ''michaels'', ''undertaker'', ''championship'', ''pinfall'', ''wrestlers'', ''heavyweight'', ''tag'', ''raw'', ''triple'', ''wrestling'', ''wwe'', ''matches'']'
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 5).Value = 1
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = 1

$ws.Cells.Item(29, 2).Value = '[''bach'', ''text'', ''conscience'', ''jesus'', ''movements'', ''movement'', ''soprano'', ''gospel'', ''hebrew'', ''aria'', ''alto'', ''manuscripts'', ''leipzig'', ''manuscript'', ''tenor'']'
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 5).Value = 1
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 7).Value = 1

$ws.Cells.Item(30, 2).Value = '[''wine'', ''bacon'', ''chicken'', ''cheese'', ''recipes'', ''sandwich'', ''fried'', ''dish'', ''dishes'', ''cuisine'', ''ingredients'', ''cooking'', ''beef'', ''food'', ''meat'']'
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(30, 5).Value = 1
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 7).Value = 1

$ws.Cells.Item(31, 2).Value = '[''coaster'', ''ride'', ''roller'', ''riders'', ''train'', ''coasters'', ''flags'', ''park'', ''lift'', ''cedar'', ''drop'', ''brake'', ''trains'', ''steel'', ''hill'']'
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(31, 5).Value = 1
$ws.Cells.Item(31, 6).Value = 1
$ws.Cells.Item(31, 7).Value = 1

$ws.Cells.Item(32, 2).Value = '[''yard'', ''yards'', ''tech'', ''touchdown'', ''bowl'', ''quarter'', ''michigan'', ''offense'', ''pass'', ''conference'', ''quarterback'', ''rushing'', ''ball'', ''football'', ''coach'']'
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(32, 5).Value = 1
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(32, 7).Value = 1

$ws.Cells.Item(33, 2).Value = '[''athletes'', ''olympics'', ''olympic'', ''beijing'', ''heat'', ''seconds'', ''athlete'', ''round'', ''games'', ''paralympic'', ''summer'', ''medal'', ''medals'', ''event'', ''events'']'
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(33, 5).Value = 1
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 7).Value = 1

$ws.Cells.Item(34, 2).Value = '[''comics'', ''fiction'', ''stories'', ''gay'', ''magazine'', ''pulp'', ''comic'', ''story'', ''science'', ''lgbt'', ''adventures'', ''issue'', ''cartoonist'', ''belgian'', ''published'']'
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(34, 5).Value = 1
$ws.Cells.Item(34, 6).Value = 1
$ws.Cells.Item(34, 7).Value = 1

$ws.Cells.Item(35, 2).Value = '[''formula'', ''function'', ''matrix'', ''linear'', ''cylinders'', ''constant'', ''functions'', ''cylinder'', ''defined'', ''filter'', ''frequency'', ''voltage'', ''gas'', ''mass'', ''derivative'']'
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(35, 5).Value = 1
$ws.Cells.Item(35, 6).Value = 1
$ws.Cells.Item(35, 7).Value = 1

$ws.Cells.Item(36, 2).Value = '[''management'', ''twitter'', ''investment'', ''bank'', ''billion'', ''equity'', ''design'', ''company'', ''banking'', ''firm'', ''users'', ''business'', ''corporate'', ''merger'', ''assets'']'
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(36, 5).Value = 1
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 1

$ws.Cells.Item(37, 2).Value = '[''fischer'', ''chess'', ''fight'', ''tournament'', ''round'', ''boxing'', ''champion'', ''hughes'', ''bout'', ''trinidad'', ''robinson'', ''punches'', ''black'', ''knockout'', ''unanimous'']'
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(37, 5).Value = 1
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(37, 7).Value = 1

$ws.Cells.Item(38, 2).Value = '[''amendment'', ''court'', ''constitution'', ''clause'', ''justices'', ''shall'', ''rights'', ''congress'', ''declaration'', ''courts'', ''law'', ''states'', ''defendant'', ''justice'', ''jurisdiction'']'
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 6).Value = 1
$ws.Cells.Item(38, 7).Value = 1

$ws.Cells.Item(39, 2).Value = '[''plants'', ''plant'', ''botanical'', ''camouflage'', ''animals'', ''gardens'', ''animal'', ''garden'', ''tea'', ''species'', ''predators'', ''coloration'', ''predator'', ''organisms'', ''mimic'']'
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(39, 4).Value = 1
$ws.Cells.Item(39, 5).Value = 1
$ws.Cells.Item(39, 6).Value = 1
$ws.Cells.Item(39, 7).Value = 1

$ws.Cells.Item(40, 2).Value = '[''phillies'', ''inning'', ''yankees'', ''dodgers'', ''teams'', ''breaker'', ''giants'', ''mlb'', ''tie'', ''yankee'', ''postseason'', ''game'', ''run'', ''pitcher'', ''innings'']'
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(40, 5).Value = 1
$ws.Cells.Item(40, 6).Value = 1
$ws.Cells.Item(40, 7).Value = 1

$ws.Cells.Item(41, 2).Value = '[''bond'', ''film'', ''casino'', ''agent'', ''films'', ''kill'', ''sequence'', ''filming'', ''majesty'', ''secret'', ''licence'', ''spy'', ''scenes'', ''villain'', ''chase'']'
$ws.Cells.Item(41, 3).Value = 1
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(41, 5).Value = 1
$ws.Cells.Item(41, 6).Value = 1
$ws.Cells.Item(41, 7).Value = 1

$ws.Cells.Item(42, 2).Value = '[''baptism'', ''congregation'', ''christ'', ''churches'', ''church'', ''spirit'', ''congregations'', ''holy'', ''assemblies'', ''god'', ''worship'', ''jewish'', ''faith'', ''israel'', ''christian'']'
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(42, 6).Value = 1
$ws.Cells.Item(42, 7).Value = 1

$ws.Cells.Item(43, 2).Value = '[''painting'', ''paintings'', ''symphony'', ''art'', ''works'', ''swan'', ''white'', ''abstract'', ''artist'', ''canvas'', ''dots'', ''girl'', ''prelude'', ''okay'', ''work'']'
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(43, 4).Value = 1
$ws.Cells.Item(43, 5).Value = 1
$ws.Cells.Item(43, 6).Value = 1
$ws.Cells.Item(43, 7).Value = 1

$ws.Cells.Item(44, 2).Value = '[''contest'', ''broadcaster'', ''semi'', ''countries'', ''jury'', ''final'', ''participating'', ''greece'', ''voting'', ''host'', ''idol'', ''entry'', ''philippine'', ''song'', ''country'']'
$ws.Cells.Item(44, 3).Value = 1
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(44, 5).Value = 1
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(44, 7).Value = 1

$ws.Cells.Item(45, 2).Value = '[''motorway'', ''croatia'', ''adriatic'', ''croatian'', ''toll'', ''traffic'', ''interchanges'', ''kilometre'', ''route'', ''interchange'', ''kilometres'', ''section'', ''areas'', ''rest'', ''yugoslavia'']'
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(45, 5).Value = 1
$ws.Cells.Item(45, 6).Value = 1
$ws.Cells.Item(45, 7).Value = 1

$ws.Cells.Item(46, 2).Value = '[''singapore'', ''law'', ''judicial'', ''constitution'', ''parliament'', ''article'', ''courts'', ''minister'', ''court'', ''constitutional'', ''detention'', ''persons'', ''public'', ''act'', ''mps'']'
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(46, 5).Value = 1
$ws.Cells.Item(46, 6).Value = 1
$ws.Cells.Item(46, 7).Value = 1

$ws.Cells.Item(47, 2).Value = '[''lighthouse'', ''light'', ''keeper'', ''tower'', ''lens'', ''keepers'', ''concrete'', ''lamp'', ''cottages'', ''constructed'', ''lamps'', ''connecticut'', ''installed'', ''iron'', ''restoration'']'
$ws.Cells.Item(47, 3).Value = 1
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(47, 5).Value = 1
$ws.Cells.Item(47, 6).Value = 1
$ws.Cells.Item(47, 7).Value = 1

$ws.Cells.Item(48, 2).Value = '[''clark'', ''superman'', ''oliver'', ''finale'', ''season'', ''comic'', ''character'', ''relationship'', ''whitney'', ''metropolis'', ''believes'', ''discovers'', ''secret'', ''martha'', ''series'']'
$ws.Cells.Item(48, 3).Value = 1
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(48, 5).Value = 1
$ws.Cells.Item(48, 6).Value = 1
$ws.Cells.Item(48, 7).Value = 1

$ws.Cells.Item(49, 2).Value = '[''euro'', ''coins'', ''currency'', ''note'', ''notes'', ''denominations'', ''dollar'', ''value'', ''stripe'', ''tender'', ''silver'', ''thread'', ''ink'', ''issued'', ''signature'']'
$ws.Cells.Item(49, 3).Value = 1
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(49, 5).Value = 1
$ws.Cells.Item(49, 6).Value = 1
$ws.Cells.Item(49, 7).Value = 1

$ws.Cells.Item(50, 2).Value = '[''grammy'', ''hawaiian'', ''nominees'', ''category'', ''awards'', ''rap'', ''award'', ''presented'', ''categories'', ''recipients'', ''academy'', ''disco'', ''honor'', ''proficiency'', ''best'']'
$ws.Cells.Item(50, 3).Value = 1
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(50, 5).Value = 1
$ws.Cells.Item(50, 6).Value = 1
$ws.Cells.Item(50, 7).Value = 1
